# Complete user2 QC Excel sheet
# - MAIN_CONTROLLER: row 6 "Y" -> "N"; add new row 7 (QC row for
#   CPC_1stTouchPoint_Approval) mirroring the layout used by earlier rows.
# - DATASHEET: row 7 (duplicate of row 6) removed, shifting the old row 8
#   (the real "FOS7_AfterPostSanction" QC record) up to row 7.
# - Active tab / selection state updated to reflect DATASHEET being the
#   sheet left active after the edit.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("MAIN_CONTROLLER")
$ws2 = $wb.Worksheets.Item("DATASHEET")

# --- MAIN_CONTROLLER ---------------------------------------------------

# Row 6: Process result flips from Y to N.
$ws1.Range("B6").Value2 = "N"

# New row 7: same shape as the other QC rows (A: quote-prefixed Si_No,
# B: Y, C: local, D/E: the new FOS/datasheet name), written before the
# formatting is pasted over it so the text values are not clobbered.
$ws1.Range("A7").Value2 = "'5"
$ws1.Range("B7").Value2 = "Y"
$ws1.Range("C7").Value2 = "local"
$ws1.Range("D7").Value2 = "CPC_1stTouchPoint_Approval"
$ws1.Range("E7").Value2 = "CPC_1stTouchPoint_Approval"

# D7:E7 pick up the highlighted-cell style already used on DATASHEET!C6.
$ws2.Range("C6").Copy() | Out-Null
$ws1.Range("D7:E7").PasteSpecial(-4122) | Out-Null

# --- DATASHEET -----------------------------------------------------------

# Row 7 was a duplicate of row 6; remove it so the real row 8 record
# (FOS7_AfterPostSanction) shifts up into row 7.
$ws2.Rows("7:7").Delete()

# --- View / selection state ---------------------------------------------

$ws1.Range("B7").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("E22").Select() | Out-Null
